# The deck currently has the "Integral" design applied (custom green/gold
# colour scheme) while its notes master still carries the stock "Office
# Theme" palette. This commit reverts the presentation's applied theme
# colours back to the default Office Theme values (the swap recorded in
# ppt/theme/theme1.xml <-> ppt/theme/theme2.xml).
#
# PowerPoint's object model exposes the live (slide-master-facing) theme
# through Slide.ThemeColorScheme - a 12-slot palette in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# Writing to it rewrites the <a:clrScheme> of the theme part shared by the
# slide master (and, transitively, every slide/layout in the deck).

function Convert-RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Stock "Office Theme" colour scheme (the target palette).
$officeThemeColors = @(
    (Convert-RGB 0x00 0x00 0x00),  # 1  dk1      000000
    (Convert-RGB 0xFF 0xFF 0xFF),  # 2  lt1      FFFFFF
    (Convert-RGB 0x44 0x54 0x6A),  # 3  dk2      44546A
    (Convert-RGB 0xE7 0xE6 0xE6),  # 4  lt2      E7E6E6
    (Convert-RGB 0x5B 0x9B 0xD5),  # 5  accent1  5B9BD5
    (Convert-RGB 0xED 0x7D 0x31),  # 6  accent2  ED7D31
    (Convert-RGB 0xA5 0xA5 0xA5),  # 7  accent3  A5A5A5
    (Convert-RGB 0xFF 0xC0 0x00),  # 8  accent4  FFC000
    (Convert-RGB 0x44 0x72 0xC4),  # 9  accent5  4472C4
    (Convert-RGB 0x70 0xAD 0x47),  # 10 accent6  70AD47
    (Convert-RGB 0x05 0x63 0xC1),  # 11 hlink    0563C1
    (Convert-RGB 0x95 0x4F 0x72)   # 12 folHlink 954F72
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColors[$i - 1]
}
